$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values per the repulled/recalculated data
$ws.Range("F5").Value = -1
$ws.Range("F6").Value = -5
$ws.Range("F8").Value = -1
$ws.Range("F9").Value = 4
$ws.Range("F10").Value = -1
$ws.Range("F11").Value = 1
$ws.Range("F12").Value = 3
$ws.Range("F13").Value = -3
$ws.Range("F15").Value = 2
$ws.Range("F16").Value = 1
$ws.Range("F17").Value = -2
$ws.Range("F18").Value = 3

$wb.Save()
